$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing xpath values
$ws.Range("B25").Value = "(//button[text()=' Open '])[99]"
$ws.Range("B101").Value = "//div[@class='mt-2']/div/table/tbody/tr[2]/td[2]/div/a/div"

# Add new rows for the Product Register module
$ws.Range("A117").Value = "Filter.nonselect.result"
$ws.Range("B117").Value = "//div[@class='table-border-dark table-responsive-sm']/table/tbody/tr/td/div/h4"

$ws.Range("A118").Value = "Product.values"
$ws.Range("B118").Value = "Study Notes"

$ws.Range("A119").Value = "search.field.value"
$ws.Range("B119").Value = "3/SL-24"

$ws.Range("A120").Value = "filter.button"
$ws.Range("B120").Value = "//main[@class='mb-5']/section/div[2]/section/div[3]/div/button"

# Move selection to mirror where the author's cursor ended up after adding the rows
$ws.Range("A122").Select()

